$d = $word.ActiveDocument

# Go to the very end of the document (after "Next, metadata filters.")
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range

# Three blank paragraphs
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()

# A fourth new blank paragraph - this will become the "date" subtitle line
$tailRange.InsertParagraphAfter()
$datePara = $d.Paragraphs.Last
$dateRange = $datePara.Range

# Create the final log-entry paragraph now, while the date paragraph is
# still plain/Normal, so the new paragraph doesn't inherit the Subtitle
# style that gets applied to the date paragraph afterwards.
$dateRange.InsertParagraphAfter()
$logPara = $d.Paragraphs.Last
$logRange = $logPara.Range
$logRange.InsertAfter("Ok I added thumbnail for both pictures and videos, and fixed layout/aspect ratio so it looks better.")

# Now fill in and style the date paragraph: "11th March 2024"
$dateRange.InsertAfter("11th March 2024")
$datePara.Style = "Subtitle"

# Make the "th" ordinal suffix superscript, like the "10th" entry above it
$supRange = $datePara.Range.Duplicate
$supRange.Find.Execute("th", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$supRange.Font.Superscript = $true

Write-Output "done"
